$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 0.651114
$ws.Range("I2").Value = 0.8572432933444277
$ws.Range("J2").Value = 0.8572432933444277
$ws.Range("Q2").Value = 0.079642745214
$ws.Range("R2").Value = 0.7167847069259999
$ws.Range("S2").Value = 0.8572432933444277
$ws.Range("T2").Value = 0.8572432933444277

# Row 3 updates
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03614333333333333
$ws.Range("H3").Value = 0.10843
$ws.Range("I3").Value = 0.1427567066555723
$ws.Range("J3").Value = 0.1427567066555723
$ws.Range("Q3").Value = 0.01326290459666667
$ws.Range("R3").Value = 0.11936614137
$ws.Range("S3").Value = 0.1427567066555723
$ws.Range("T3").Value = 0.1427567066555723
